$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), matching the style of the existing header (H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I and J columns (I0 / IF) for rows 2-70
$data = @(
  @(2, 9, 9),
  @(3, 8, 9),
  @(4, 8, 8),
  @(5, 7, 8),
  @(6, 6, 6),
  @(7, 8, 8),
  @(8, 8, 8),
  @(9, 8, 8),
  @(10, 8, 8),
  @(11, 8, 8),
  @(12, 8, 8),
  @(13, 8, 8),
  @(14, 7, 8),
  @(15, 8, 8),
  @(16, 11, 12),
  @(17, 8, 8),
  @(18, 8, 8),
  @(19, 9, 9),
  @(20, 8, 8),
  @(21, 8, 8),
  @(22, 7, 7),
  @(23, 6, 7),
  @(24, 7, 7),
  @(25, 6, 6),
  @(26, 6, 6),
  @(27, 6, 6),
  @(28, 4, 4),
  @(29, 8, 8),
  @(30, 9, 9),
  @(31, 8, 8),
  @(32, 6, 6),
  @(33, 6, 6),
  @(34, 6, 6),
  @(35, 9, 9),
  @(36, 6, 6),
  @(37, 6, 6),
  @(38, 5, 6),
  @(39, 6, 6),
  @(40, 6, 7),
  @(41, 9, 9),
  @(42, 7, 7),
  @(43, 8, 8),
  @(44, 8, 8),
  @(45, 8, 8),
  @(46, 9, 9),
  @(47, 9, 9),
  @(48, 9, 9),
  @(49, 9, 9),
  @(50, 7, 7),
  @(51, 8, 8),
  @(52, 9, 9),
  @(53, 7, 8),
  @(54, 9, 9),
  @(55, 8, 9),
  @(56, 9, 9),
  @(57, 9, 9),
  @(58, 9, 9),
  @(59, 6, 7),
  @(60, 8, 8),
  @(61, 6, 6),
  @(62, 7, 8),
  @(63, 7, 7),
  @(64, 6, 7),
  @(65, 9, 9),
  @(66, 6, 7),
  @(67, 7, 7),
  @(68, 4, 4),
  @(69, 6, 6),
  @(70, 9, 9)
)

foreach ($entry in $data) {
  $r = $entry[0]
  $iVal = $entry[1]
  $jVal = $entry[2]
  $ws.Cells.Item($r, 9).Value = $iVal
  $ws.Cells.Item($r, 10).Value = $jVal
}

"Completed: added I0/IF columns"